$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 23.86000000000029
$ws.Range("H2").Value = [double]"5.163828021512356e-16"
$ws.Range("K2").Value = 35.1035621680703
$ws.Range("L2").Value = "[28.539370974147076, 41.66775336199352]"
$ws.Range("O2").Value = 1.566079220708425
$ws.Range("P2").Value = "[1.3522370781217328, 1.7799213632951174]"
$ws.Range("S2").Value = 55.72682124626629
$ws.Range("T2").Value = "[51.32811479012665, 60.12552770240593]"
$ws.Range("W2").Value = 17.91291291291313
$ws.Range("X2").Value = 17.10086086086107
$ws.Range("Y2").Value = 18.72496496496519

# Row 3 updates
$ws.Range("E3").Value = 24.05000000000032
$ws.Range("G3").Value = [double]"2.220446049250313e-13"
$ws.Range("H3").Value = [double]"7.757740968750056e-13"
$ws.Range("K3").Value = 39.80500862369955
$ws.Range("L3").Value = "[26.799188434170482, 52.81082881322862]"
$ws.Range("M3").Value = [double]"8.136850970785758e-09"
$ws.Range("N3").Value = [double]"8.136850970785758e-09"
$ws.Range("O3").Value = 2.182447749340658
$ws.Range("P3").Value = "[1.8553950606786582, 2.5095004380026573]"
$ws.Range("S3").Value = 59.78008777572396
$ws.Range("T3").Value = "[53.059098085979386, 66.50107746546854]"
$ws.Range("W3").Value = 15.6962962962965
$ws.Range("X3").Value = 14.44444444444464
$ws.Range("Y3").Value = 16.94814814814837
